$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "patient"
$ws.Range("B1").Value = "mood_pre"
$ws.Range("C1").Value = "mood_post"

$ws.Range("A1").Select()
